$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column U: "Дата вывода из эксплуатации" (archive/decommission date) ---

# Header cell U1 — copy formatting from T1 (same header style), then set text.
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)
$ws.Range("U1").Value = "Дата вывода из эксплуатации"

# Data cells U2:U11 — copy formatting from T2:T11 (same data-row style), then
# set text. These are a copy of the meter-status-column styling.
$ws.Range("T2:T11").Copy()
$ws.Range("U2:U11").PasteSpecial(-4122)

# Column width to roughly match the authored sheet (~24.17 OOXML width units).
$ws.Columns.Item(21).ColumnWidth = 23.25

# --- Example data values ---

# Row 2 ("Автоматический" example row): mark as decommissioned + give a date.
$ws.Range("T2").Value = "Да"
$ws.Range("U2").Value = "2022-01-23"

# Row 3: explicitly mark as not decommissioned.
$ws.Range("T3").Value = "Нет"

$excel.CutCopyMode = 0
